{"js": "// \"Version 1.\" -> \"Version 2.\" while reproducing the exact run-split\n// structure from the target revision:\n//   - \"Version\" splits into two runs: \"Versi\" + \"on\" (spellcheck markers\n//     stay wrapped tightly around both pieces)\n//   - \" 1.\" loses its trailing \".\" and becomes \" 2\"\n//   - a brand new run containing \".\" is appended *after* the existing\n//     _GoBack bookmark\nconst doc = context.document;\nconst body = doc.body;\n\n// Step 1: capture the bookmark's (zero-width) position while \" 1.\" is\n// still intact, and insert the new trailing \".\" run right after it.\n// Doing this first means the later text replacements don't shift where\n// the bookmark range resolves to.\nconst bookmarkRange = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nbookmarkRange.insertText(\".\", Word.InsertLocation.after);\nawait context.sync();\n\n// Step 2: replace the original \" 1.\" (the occurrence that still precedes\n// the bookmark) with \" 2\" - drops the period, flips the digit.\nconst versionNumber = body.search(\" 1.\", { matchCase: true });\nawait context.sync();\nversionNumber.items[0].insertText(\" 2\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 3: split \"Version\" into \"Versi\" + \"on\" runs. Replacing only the\n// \"on\" suffix (so the edited range touches the spellEnd proofErr boundary\n// but not spellStart) keeps the proof-error markers anchored around the\n// whole word instead of being pushed out of position.\nconst suffix = body.search(\"on\", { matchCase: true });\nawait context.sync();\nconst splitRunsOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>on</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nsuffix.items[0].insertOoxml(splitRunsOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"Version 1.\" -> \"Version 2.\" while reproducing the exact run-split\n# structure from the target revision:\n#   - \"Version\" splits into two runs: \"Versi\" + \"on\" (spellcheck markers\n#     stay wrapped tightly around both pieces)\n#   - \" 1.\" loses its trailing \".\" and becomes \" 2\"\n#   - a brand new run containing \".\" is appended *after* the existing\n#     _GoBack bookmark\n$d = $word.ActiveDocument\n\n# Step 1: insert the new trailing \".\" run right after the existing\n# _GoBack bookmark while \" 1.\" is still intact, so later edits don't\n# shift where the bookmark resolves to.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bmRange = $d.Range($bm.Start, $bm.End)\n$bmRange.InsertAfter(\".\")\n\n# Step 2: replace the original \" 1.\" (the occurrence that still precedes\n# the bookmark) with \" 2\" - drops the period, flips the digit.\n$find = $d.Content\n$find.Find.Execute(\" 1.\", $false, $false, $false, $false, $false, $true, 1, $false, \" 2\", 2) | Out-Null\n\n# Step 3: split \"Version\" into \"Versi\" + \"on\" runs. A temporary bookmark\n# dropped at the split point (character offset 5) forces Word to break\n# the run there; deleting the bookmark right after removes the marker\n# but leaves the two runs in place with no leftover formatting.\n$splitPoint = $d.Range(5, 5)\n$d.Bookmarks.Add(\"TempSplit\", $splitPoint) | Out-Null\n$d.Bookmarks.Item(\"TempSplit\").Delete()\n"}
